# The Google-sheet source re-generated this .ttl: a new "owl" PREFIX row was
# inserted into the "Prefixes" block (pushing every following row down by
# one), and the dct:modified^^xsd:datetime timestamp was bumped to reflect
# the regeneration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 (old rows 13-30 shift down to 14-31).
$ws.Rows.Item(13).Insert()

# New row 13: PREFIX | owl | http://www.w3.org/2002/07/owl#
$ws.Range("A13").Value = "PREFIX"
$ws.Range("B13").Value = "owl"
$ws.Range("C13").Value = "http://www.w3.org/2002/07/owl#"

# dct:modified^^xsd:datetime row is now row 22 (was row 21 pre-insert);
# bump its value to the regeneration timestamp.
$ws.Range("B22").Value = "2023-09-13T15:17:21+00:00"
